$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 384.2857
$ws.Range("I2").Value = 398.33334
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 398.33334
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -285.33334
$ws.Range("N2").Value = -526

$ws.Range("H29").Value = 900
$ws.Range("I29").Value = 900
$ws.Range("K29").Value = 2700
$ws.Range("M29").Value = -2419

$ws.Range("H33").Value = 129.57143
$ws.Range("I33").Value = 125.75
$ws.Range("K33").Value = 125.75
$ws.Range("M33").Value = 103.25

$ws.Range("H38").Value = 1380.9286
$ws.Range("I38").Value = 137.5
$ws.Range("J38").Value = 3038.8333
$ws.Range("K38").Value = 412.5
$ws.Range("L38").Value = 9116.499899999999
$ws.Range("M38").Value = -40.5
$ws.Range("N38").Value = -9860.499899999999

$ws.Range("H74").Value = 3500
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2064
$ws.Range("N74").Value = -5872

$ws.Range("H77").Value = 3500
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -29360

$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H129").Value = 1187.4242
$ws.Range("J129").Value = 1231.9678
$ws.Range("L129").Value = 3695.9034
$ws.Range("N129").Value = -13695.9034

$ws.Range("H137").Value = 34484056
$ws.Range("I137").Value = 52632460
$ws.Range("J137").Value = 2086.5
$ws.Range("K137").Value = 157897380
$ws.Range("L137").Value = 6259.5
$ws.Range("M137").Value = -157894830
$ws.Range("N137").Value = -11359.5

$ws.Range("H138").Value = 4164656.2
$ws.Range("I138").Value = 1034979.9
$ws.Range("J138").Value = 7094566
$ws.Range("K138").Value = 3104939.7
$ws.Range("L138").Value = 21283698
$ws.Range("M138").Value = -3099799.7
$ws.Range("N138").Value = -21293978

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -201

$ws.Range("H28").Value = 108194.2
$ws.Range("I28").Value = 126992.75
$ws.Range("J28").Value = 33000
$ws.Range("K28").Value = 126992.75
$ws.Range("L28").Value = 33000
$ws.Range("M28").Value = -126800.75
$ws.Range("N28").Value = -33384

$ws.Range("H32").Value = 13870.528
$ws.Range("I32").Value = 2518.061
$ws.Range("J32").Value = 146856.58
$ws.Range("K32").Value = 2518.061
$ws.Range("L32").Value = 146856.58
$ws.Range("M32").Value = -2231.061
$ws.Range("N32").Value = -147430.58

$ws.Range("H74").Value = 5409.2163
$ws.Range("I74").Value = 1986.8077
$ws.Range("J74").Value = 13498.546
$ws.Range("K74").Value = 1986.8077
$ws.Range("L74").Value = 13498.546
$ws.Range("M74").Value = -1112.8077
$ws.Range("N74").Value = -15246.546

$ws.Range("H77").Value = 5409.2163
$ws.Range("I77").Value = 1986.8077
$ws.Range("J77").Value = 13498.546
$ws.Range("K77").Value = 9934.038500000001
$ws.Range("L77").Value = 67492.73
$ws.Range("M77").Value = -5566.038500000001
$ws.Range("N77").Value = -76228.73

$ws.Range("H97").Value = 6188.222
$ws.Range("I97").Value = 6867.375
$ws.Range("J97").Value = 755
$ws.Range("K97").Value = 6867.375
$ws.Range("L97").Value = 755
$ws.Range("M97").Value = -6371.375
$ws.Range("N97").Value = -1747

$ws.Range("H99").Value = 108194.2
$ws.Range("I99").Value = 126992.75
$ws.Range("J99").Value = 33000
$ws.Range("K99").Value = 126992.75
$ws.Range("L99").Value = 33000
$ws.Range("M99").Value = -123997.75
$ws.Range("N99").Value = -38990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8957.071
$ws.Range("I86").Value = 2059.2
$ws.Range("J86").Value = 26201.75
$ws.Range("K86").Value = 2059.2
$ws.Range("L86").Value = 26201.75
$ws.Range("M86").Value = -936.1999999999998
$ws.Range("N86").Value = -28447.75

$ws.Range("H89").Value = 8957.071
$ws.Range("I89").Value = 2059.2
$ws.Range("J89").Value = 26201.75
$ws.Range("K89").Value = 10296
$ws.Range("L89").Value = 131008.75
$ws.Range("M89").Value = -4680
$ws.Range("N89").Value = -142240.75

$ws.Range("H99").Value = 2090.818
$ws.Range("I99").Value = 2090.818
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2090.818
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -592.8180000000002
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 255134.28
$ws.Range("I105").Value = 6104.7407
$ws.Range("J105").Value = 772349.4399999999
$ws.Range("K105").Value = 6104.7407
$ws.Range("L105").Value = 772349.4399999999
$ws.Range("M105").Value = -4357.7407
$ws.Range("N105").Value = -775843.4399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 39878.31
$ws.Range("I16").Value = 50991.8
$ws.Range("J16").Value = 2833.3333
$ws.Range("K16").Value = 50991.8
$ws.Range("L16").Value = 2833.3333
$ws.Range("M16").Value = -50704.8
$ws.Range("N16").Value = -3407.3333

$ws.Range("H31").Value = 1806.18
$ws.Range("I31").Value = 1002.8108
$ws.Range("K31").Value = 1002.8108
$ws.Range("M31").Value = -707.8108

$ws.Range("H34").Value = 1806.18
$ws.Range("I34").Value = 1002.8108
$ws.Range("K34").Value = 1002.8108
$ws.Range("M34").Value = -800.8108

$ws.Range("H58").Value = 1325.44
$ws.Range("I58").Value = 919.6053000000001
$ws.Range("J58").Value = 2610.5833
$ws.Range("K58").Value = 919.6053000000001
$ws.Range("L58").Value = 2610.5833
$ws.Range("M58").Value = -716.6053000000001
$ws.Range("N58").Value = -3016.5833

$ws.Range("H107").Value = 414.1875
$ws.Range("I107").Value = 235.5
$ws.Range("J107").Value = 592.875
$ws.Range("K107").Value = 235.5
$ws.Range("L107").Value = 592.875
$ws.Range("M107").Value = 1684.5
$ws.Range("N107").Value = -4432.875

$ws.Range("H109").Value = 26719.334
$ws.Range("I109").Value = 20258
$ws.Range("J109").Value = 29950
$ws.Range("K109").Value = 20258
$ws.Range("L109").Value = 29950
$ws.Range("M109").Value = -19218
$ws.Range("N109").Value = -32030

$ws.Range("H113").Value = 39878.31
$ws.Range("I113").Value = 50991.8
$ws.Range("J113").Value = 2833.3333
$ws.Range("K113").Value = 50991.8
$ws.Range("L113").Value = 2833.3333
$ws.Range("M113").Value = -48821.8
$ws.Range("N113").Value = -7173.3333

$ws.Range("H122").Value = 1608
$ws.Range("I122").Value = 1606.5
$ws.Range("K122").Value = 4819.5
$ws.Range("M122").Value = -2369.5

$ws.Range("H136").Value = 1325.44
$ws.Range("I136").Value = 919.6053000000001
$ws.Range("J136").Value = 2610.5833
$ws.Range("K136").Value = 2758.8159
$ws.Range("L136").Value = 7831.749899999999
$ws.Range("M136").Value = -208.8159000000001
$ws.Range("N136").Value = -12931.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 41.75
$ws.Range("I14").Value = 41.75
$ws.Range("K14").Value = 125.25
$ws.Range("M14").Value = 47.75

$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H74").Value = 12000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 12000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 36000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -38122

$ws.Range("H77").Value = 12000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 12000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 108000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -118608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7215.476
$ws.Range("I70").Value = 7786.5386
$ws.Range("J70").Value = 6287.5
$ws.Range("K70").Value = 7786.5386
$ws.Range("L70").Value = 6287.5
$ws.Range("M70").Value = -7516.5386
$ws.Range("N70").Value = -6827.5

$ws.Range("H73").Value = 7215.476
$ws.Range("I73").Value = 7786.5386
$ws.Range("J73").Value = 6287.5
$ws.Range("K73").Value = 7786.5386
$ws.Range("L73").Value = 6287.5
$ws.Range("M73").Value = -6850.5386
$ws.Range("N73").Value = -8159.5

$ws.Range("H132").Value = 3912.6428
$ws.Range("I132").Value = 3742.125
$ws.Range("J132").Value = 4935.75
$ws.Range("K132").Value = 11226.375
$ws.Range("L132").Value = 14807.25
$ws.Range("M132").Value = -8696.375
$ws.Range("N132").Value = -19867.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1958.0834
$ws.Range("I68").Value = 1437.125
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1437.125
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -688.125
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 1958.0834
$ws.Range("I71").Value = 1437.125
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 7185.625
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -3441.625
$ws.Range("N71").Value = -22488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 34195
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 34195
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 34195
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -35703

$ws.Range("H74").Value = 10700.091
$ws.Range("J74").Value = 10329.125
$ws.Range("L74").Value = 10329.125
$ws.Range("N74").Value = -12201.125

$ws.Range("H77").Value = 10700.091
$ws.Range("J77").Value = 10329.125
$ws.Range("L77").Value = 30987.375
$ws.Range("N77").Value = -40347.375

$ws.Range("H109").Value = 32125.666
$ws.Range("J109").Value = 32125.666
$ws.Range("L109").Value = 32125.666
$ws.Range("N109").Value = -34899.666
